$wb = $excel.ActiveWorkbook

# --- Step1_Data: update signal-value distribution cells (rows 4,5,8,11) ---
$ws1 = $wb.Worksheets.Item("Step1_Data")
# row 4
$ws1.Cells.Item(4, 13).Value = 0
$ws1.Cells.Item(4, 14).Value = 0
$ws1.Cells.Item(4, 15).Value = 0
$ws1.Cells.Item(4, 16).Value = 0
$ws1.Cells.Item(4, 17).Value = 0
$ws1.Cells.Item(4, 18).Value = 0
$ws1.Cells.Item(4, 19).Value = 0
$ws1.Cells.Item(4, 20).Value = 0
$ws1.Cells.Item(4, 21).Value = 0
$ws1.Cells.Item(4, 22).Value = 0
$ws1.Cells.Item(4, 23).Value = 0
$ws1.Cells.Item(4, 24).Value = 0
$ws1.Cells.Item(4, 25).Value = 0
$ws1.Cells.Item(4, 26).Value = 0
$ws1.Cells.Item(4, 27).Value = 0
$ws1.Cells.Item(4, 28).Value = 0
$ws1.Cells.Item(4, 29).Value = 0
$ws1.Cells.Item(4, 30).Value = 0
$ws1.Cells.Item(4, 31).Value = 0
$ws1.Cells.Item(4, 32).Value = 0
$ws1.Cells.Item(4, 33).Value = 0
$ws1.Cells.Item(4, 34).Value = 0
$ws1.Cells.Item(4, 35).Value = 0
$ws1.Cells.Item(4, 36).Value = 0
$ws1.Cells.Item(4, 37).Value = 0
$ws1.Cells.Item(4, 38).Value = 0
$ws1.Cells.Item(4, 39).Value = 0
$ws1.Cells.Item(4, 40).Value = 0
$ws1.Cells.Item(4, 41).Value = 0
$ws1.Cells.Item(4, 42).Value = 0
$ws1.Cells.Item(4, 43).Value = 0
$ws1.Cells.Item(4, 44).Value = 0
$ws1.Cells.Item(4, 93).Value = 0.1802040037029629
$ws1.Cells.Item(4, 94).Value = 0.00029303490855601
$ws1.Cells.Item(4, 95).Value = 0.216484675194083
$ws1.Cells.Item(4, 96).Value = 0.004499580771517468
$ws1.Cells.Item(4, 97).Value = 0.01968345639726003
$ws1.Cells.Item(4, 98).Value = 0.05930032475074898
$ws1.Cells.Item(4, 99).Value = 0.006823657184081323
$ws1.Cells.Item(4, 100).Value = 0.003393463501923142
$ws1.Cells.Item(4, 101).Value = 0.002662479784846849
$ws1.Cells.Item(4, 102).Value = 0.03004135118908167
$ws1.Cells.Item(4, 103).Value = 0.03854708820318444
$ws1.Cells.Item(4, 104).Value = 0.02360236069387227
$ws1.Cells.Item(4, 105).Value = 0.07265879782828374
$ws1.Cells.Item(4, 106).Value = 0.001610092473027443
$ws1.Cells.Item(4, 107).Value = 0.1568543070916621
$ws1.Cells.Item(4, 108).Value = 0.02865204325914033
$ws1.Cells.Item(4, 109).Value = 0.07462869621042069
$ws1.Cells.Item(4, 110).Value = 0.02330665925466352
$ws1.Cells.Item(4, 111).Value = 0.01015574294098501
$ws1.Cells.Item(4, 112).Value = 0.00003446935391569198
$ws1.Cells.Item(4, 113).Value = 0.000210751759966172
$ws1.Cells.Item(4, 114).Value = 0.01285765839612315
$ws1.Cells.Item(4, 115).Value = 0.004712445374907431
$ws1.Cells.Item(4, 116).Value = 0.0009950823195708325
$ws1.Cells.Item(4, 117).Value = 0.00002519718272750968
$ws1.Cells.Item(4, 118).Value = 0.003621663115905833
$ws1.Cells.Item(4, 119).Value = 0.0005201620772011712
$ws1.Cells.Item(4, 120).Value = 0.01214583022469258
$ws1.Cells.Item(4, 121).Value = 0.01024172296427424
$ws1.Cells.Item(4, 122).Value = 0.0003626380280612266
$ws1.Cells.Item(4, 123).Value = 0.0005305249318414142
$ws1.Cells.Item(4, 124).Value = 0.0003400389305116147
# row 5
$ws1.Cells.Item(5, 26).Value = 0
$ws1.Cells.Item(5, 27).Value = 0
$ws1.Cells.Item(5, 28).Value = 0
$ws1.Cells.Item(5, 29).Value = 0
$ws1.Cells.Item(5, 30).Value = 0
$ws1.Cells.Item(5, 31).Value = 0
$ws1.Cells.Item(5, 32).Value = 0
$ws1.Cells.Item(5, 33).Value = 0
$ws1.Cells.Item(5, 34).Value = 0
$ws1.Cells.Item(5, 35).Value = 0
$ws1.Cells.Item(5, 36).Value = 0
$ws1.Cells.Item(5, 37).Value = 0
$ws1.Cells.Item(5, 38).Value = 0
$ws1.Cells.Item(5, 39).Value = 0
$ws1.Cells.Item(5, 40).Value = 0
$ws1.Cells.Item(5, 41).Value = 0
$ws1.Cells.Item(5, 42).Value = 0
$ws1.Cells.Item(5, 43).Value = 0
$ws1.Cells.Item(5, 44).Value = 0
$ws1.Cells.Item(5, 45).Value = 0
$ws1.Cells.Item(5, 46).Value = 0
$ws1.Cells.Item(5, 47).Value = 0
$ws1.Cells.Item(5, 48).Value = 0
$ws1.Cells.Item(5, 49).Value = 0
$ws1.Cells.Item(5, 50).Value = 0
$ws1.Cells.Item(5, 51).Value = 0
$ws1.Cells.Item(5, 52).Value = 0
$ws1.Cells.Item(5, 53).Value = 0
$ws1.Cells.Item(5, 54).Value = 0
$ws1.Cells.Item(5, 55).Value = 0
$ws1.Cells.Item(5, 56).Value = 0
$ws1.Cells.Item(5, 57).Value = 0
$ws1.Cells.Item(5, 93).Value = 0.0414397235972116
$ws1.Cells.Item(5, 94).Value = 0.2630806270126066
$ws1.Cells.Item(5, 95).Value = 0.2012933290688982
$ws1.Cells.Item(5, 96).Value = 0.05364909550629726
$ws1.Cells.Item(5, 97).Value = 0.001190238352764568
$ws1.Cells.Item(5, 98).Value = 0.04370335695455926
$ws1.Cells.Item(5, 99).Value = 0.02025541867063316
$ws1.Cells.Item(5, 100).Value = 0.04511906604375923
$ws1.Cells.Item(5, 101).Value = 0.0130341148681117
$ws1.Cells.Item(5, 102).Value = 0.02243130709094735
$ws1.Cells.Item(5, 103).Value = 0.04246888512389326
$ws1.Cells.Item(5, 104).Value = 0.0005586733277251922
$ws1.Cells.Item(5, 105).Value = 0.06985394093788058
$ws1.Cells.Item(5, 106).Value = 0.00097695327955818
$ws1.Cells.Item(5, 107).Value = 0.1116006711348013
$ws1.Cells.Item(5, 108).Value = 0.008923270608888102
$ws1.Cells.Item(5, 109).Value = 0.01246576125352156
$ws1.Cells.Item(5, 110).Value = 0.007083692586443485
$ws1.Cells.Item(5, 111).Value = 0.004279224405838725
$ws1.Cells.Item(5, 112).Value = 0.0005806953029109944
$ws1.Cells.Item(5, 113).Value = 0.01275590463982641
$ws1.Cells.Item(5, 114).Value = 0.002244184492132652
$ws1.Cells.Item(5, 115).Value = 0.000377997805546931
$ws1.Cells.Item(5, 116).Value = 0.0001606166517830405
$ws1.Cells.Item(5, 117).Value = 0.003325442405765627
$ws1.Cells.Item(5, 118).Value = 0.000002321932010179281
$ws1.Cells.Item(5, 119).Value = 0.001201726965841406
$ws1.Cells.Item(5, 120).Value = 0.002741142436134194
$ws1.Cells.Item(5, 121).Value = 0.000001275304341272437
$ws1.Cells.Item(5, 122).Value = 0.0006815783159268049
$ws1.Cells.Item(5, 123).Value = 0.01187032110000375
$ws1.Cells.Item(5, 124).Value = 0.0006494428234373005
# row 8
$ws1.Cells.Item(8, 17).Value = 0
$ws1.Cells.Item(8, 18).Value = 0
$ws1.Cells.Item(8, 19).Value = 0
$ws1.Cells.Item(8, 20).Value = 0
$ws1.Cells.Item(8, 21).Value = 0
$ws1.Cells.Item(8, 22).Value = 0
$ws1.Cells.Item(8, 23).Value = 0
$ws1.Cells.Item(8, 24).Value = 0
$ws1.Cells.Item(8, 25).Value = 0
$ws1.Cells.Item(8, 26).Value = 0
$ws1.Cells.Item(8, 27).Value = 0
$ws1.Cells.Item(8, 28).Value = 0
$ws1.Cells.Item(8, 29).Value = 0
$ws1.Cells.Item(8, 30).Value = 0
$ws1.Cells.Item(8, 31).Value = 0
$ws1.Cells.Item(8, 32).Value = 0
$ws1.Cells.Item(8, 33).Value = 0
$ws1.Cells.Item(8, 34).Value = 0
$ws1.Cells.Item(8, 35).Value = 0
$ws1.Cells.Item(8, 36).Value = 0
$ws1.Cells.Item(8, 37).Value = 0
$ws1.Cells.Item(8, 38).Value = 0
$ws1.Cells.Item(8, 39).Value = 0
$ws1.Cells.Item(8, 40).Value = 0
$ws1.Cells.Item(8, 41).Value = 0
$ws1.Cells.Item(8, 42).Value = 0
$ws1.Cells.Item(8, 43).Value = 0
$ws1.Cells.Item(8, 44).Value = 0
$ws1.Cells.Item(8, 45).Value = 0
$ws1.Cells.Item(8, 46).Value = 0
$ws1.Cells.Item(8, 47).Value = 0
$ws1.Cells.Item(8, 48).Value = 0
$ws1.Cells.Item(8, 92).Value = 0.14584622503048
$ws1.Cells.Item(8, 93).Value = 0.104112134841501
$ws1.Cells.Item(8, 94).Value = 0.260046837839829
$ws1.Cells.Item(8, 95).Value = 0.0419772576669443
$ws1.Cells.Item(8, 96).Value = 0.0005882564023084938
$ws1.Cells.Item(8, 97).Value = 0.02923237307107321
$ws1.Cells.Item(8, 98).Value = 0.02265675527316096
$ws1.Cells.Item(8, 99).Value = 0.01437951576595625
$ws1.Cells.Item(8, 100).Value = 0.0005930643140023797
$ws1.Cells.Item(8, 101).Value = 0.122420182502357
$ws1.Cells.Item(8, 102).Value = 0.03851413654231906
$ws1.Cells.Item(8, 103).Value = 0.01961600011603389
$ws1.Cells.Item(8, 104).Value = 0.04837384109706605
$ws1.Cells.Item(8, 105).Value = 0.0004882274881135669
$ws1.Cells.Item(8, 106).Value = 0.06157552675421839
$ws1.Cells.Item(8, 107).Value = 0.00001470227239889956
$ws1.Cells.Item(8, 108).Value = 0.01024694633098352
$ws1.Cells.Item(8, 109).Value = 0.006348813311121298
$ws1.Cells.Item(8, 110).Value = 0.00461369579762668
$ws1.Cells.Item(8, 111).Value = 0.005143365308725611
$ws1.Cells.Item(8, 112).Value = 0.01240573262854361
$ws1.Cells.Item(8, 113).Value = 0.01137202225888424
$ws1.Cells.Item(8, 114).Value = 0.003051128239303343
$ws1.Cells.Item(8, 115).Value = 0.001932049517670346
$ws1.Cells.Item(8, 116).Value = 0.0004853416986301422
$ws1.Cells.Item(8, 117).Value = 0.008707667220223965
$ws1.Cells.Item(8, 118).Value = 0.01437579068354306
$ws1.Cells.Item(8, 119).Value = 0.00851951420565939
$ws1.Cells.Item(8, 120).Value = 0.0007028352168332861
$ws1.Cells.Item(8, 121).Value = 0.000198377104136292
$ws1.Cells.Item(8, 122).Value = 0.0007820592410294135
$ws1.Cells.Item(8, 123).Value = 0.0006796242593234319
# row 11
$ws1.Cells.Item(11, 20).Value = 0
$ws1.Cells.Item(11, 21).Value = 0
$ws1.Cells.Item(11, 22).Value = 0
$ws1.Cells.Item(11, 23).Value = 0
$ws1.Cells.Item(11, 24).Value = 0
$ws1.Cells.Item(11, 25).Value = 0
$ws1.Cells.Item(11, 26).Value = 0
$ws1.Cells.Item(11, 27).Value = 0
$ws1.Cells.Item(11, 28).Value = 0
$ws1.Cells.Item(11, 29).Value = 0
$ws1.Cells.Item(11, 30).Value = 0
$ws1.Cells.Item(11, 31).Value = 0
$ws1.Cells.Item(11, 32).Value = 0
$ws1.Cells.Item(11, 33).Value = 0
$ws1.Cells.Item(11, 34).Value = 0
$ws1.Cells.Item(11, 35).Value = 0
$ws1.Cells.Item(11, 36).Value = 0
$ws1.Cells.Item(11, 37).Value = 0
$ws1.Cells.Item(11, 38).Value = 0
$ws1.Cells.Item(11, 39).Value = 0
$ws1.Cells.Item(11, 40).Value = 0
$ws1.Cells.Item(11, 41).Value = 0
$ws1.Cells.Item(11, 42).Value = 0
$ws1.Cells.Item(11, 43).Value = 0
$ws1.Cells.Item(11, 44).Value = 0
$ws1.Cells.Item(11, 45).Value = 0
$ws1.Cells.Item(11, 46).Value = 0
$ws1.Cells.Item(11, 47).Value = 0
$ws1.Cells.Item(11, 48).Value = 0
$ws1.Cells.Item(11, 49).Value = 0
$ws1.Cells.Item(11, 50).Value = 0
$ws1.Cells.Item(11, 51).Value = 0
$ws1.Cells.Item(11, 52).Value = 0.2230194153683032
$ws1.Cells.Item(11, 53).Value = 0.003029794249319589
$ws1.Cells.Item(11, 54).Value = 0.29784526246771
$ws1.Cells.Item(11, 55).Value = 0.04328837348088212
$ws1.Cells.Item(11, 56).Value = 0.007301316782442628
$ws1.Cells.Item(11, 57).Value = 0.0766752169954596
$ws1.Cells.Item(11, 58).Value = 0.02611424040570105
$ws1.Cells.Item(11, 59).Value = 0.0004013291759833063
$ws1.Cells.Item(11, 60).Value = 0.003631307137243256
$ws1.Cells.Item(11, 61).Value = 0.1130199706657033
$ws1.Cells.Item(11, 62).Value = 0.01357509079295698
$ws1.Cells.Item(11, 63).Value = 0.02821023802447835
$ws1.Cells.Item(11, 64).Value = 0.0440977578847643
$ws1.Cells.Item(11, 65).Value = 0.00004303483380139027
$ws1.Cells.Item(11, 66).Value = 0.02850124578111823
$ws1.Cells.Item(11, 67).Value = 0.006332946675936626
$ws1.Cells.Item(11, 68).Value = 0.00259913065008579
$ws1.Cells.Item(11, 69).Value = 0.000789740437905566
$ws1.Cells.Item(11, 70).Value = 0.0008570981283736624
$ws1.Cells.Item(11, 71).Value = 0.00003209637339871037
$ws1.Cells.Item(11, 72).Value = 0.003392059470871622
$ws1.Cells.Item(11, 73).Value = 0.003229992397896159
$ws1.Cells.Item(11, 74).Value = 0.0003796239287754585
$ws1.Cells.Item(11, 75).Value = 0.02916436275667109
$ws1.Cells.Item(11, 76).Value = 0.00243783855148801
$ws1.Cells.Item(11, 77).Value = 0.008109711502283515
$ws1.Cells.Item(11, 78).Value = 0.01807920619094101
$ws1.Cells.Item(11, 79).Value = 0.004984627324279639
$ws1.Cells.Item(11, 80).Value = 0.0000002320778269532873
$ws1.Cells.Item(11, 81).Value = 0.0000555743659379272
$ws1.Cells.Item(11, 82).Value = 0.002189508854612637
$ws1.Cells.Item(11, 83).Value = 0.008612656266848305

# --- Step2_Sj: update cumulative-sum cells (rows 4,5,8,11) ---
$ws2 = $wb.Worksheets.Item("Step2_Sj")
# row 4
$ws2.Cells.Item(4, 13).Value = 0
$ws2.Cells.Item(4, 14).Value = 0
$ws2.Cells.Item(4, 15).Value = 0
$ws2.Cells.Item(4, 16).Value = 0
$ws2.Cells.Item(4, 17).Value = 0
$ws2.Cells.Item(4, 18).Value = 0
$ws2.Cells.Item(4, 19).Value = 0
$ws2.Cells.Item(4, 20).Value = 0
$ws2.Cells.Item(4, 21).Value = 0
$ws2.Cells.Item(4, 22).Value = 0
$ws2.Cells.Item(4, 23).Value = 0
$ws2.Cells.Item(4, 24).Value = 0
$ws2.Cells.Item(4, 25).Value = 0
$ws2.Cells.Item(4, 26).Value = 0
$ws2.Cells.Item(4, 27).Value = 0
$ws2.Cells.Item(4, 28).Value = 0
$ws2.Cells.Item(4, 29).Value = 0
$ws2.Cells.Item(4, 30).Value = 0
$ws2.Cells.Item(4, 31).Value = 0
$ws2.Cells.Item(4, 32).Value = 0
$ws2.Cells.Item(4, 33).Value = 0
$ws2.Cells.Item(4, 34).Value = 0
$ws2.Cells.Item(4, 35).Value = 0
$ws2.Cells.Item(4, 36).Value = 0
$ws2.Cells.Item(4, 37).Value = 0
$ws2.Cells.Item(4, 38).Value = 0
$ws2.Cells.Item(4, 39).Value = 0
$ws2.Cells.Item(4, 40).Value = 0
$ws2.Cells.Item(4, 41).Value = 0
$ws2.Cells.Item(4, 42).Value = 0
$ws2.Cells.Item(4, 43).Value = 0
$ws2.Cells.Item(4, 44).Value = 0
$ws2.Cells.Item(4, 45).Value = 0
$ws2.Cells.Item(4, 46).Value = 0
$ws2.Cells.Item(4, 47).Value = 0
$ws2.Cells.Item(4, 48).Value = 0
$ws2.Cells.Item(4, 49).Value = 0
$ws2.Cells.Item(4, 50).Value = 0
$ws2.Cells.Item(4, 51).Value = 0
$ws2.Cells.Item(4, 52).Value = 0
$ws2.Cells.Item(4, 53).Value = 0
$ws2.Cells.Item(4, 54).Value = 0
$ws2.Cells.Item(4, 55).Value = 0
$ws2.Cells.Item(4, 56).Value = 0
$ws2.Cells.Item(4, 57).Value = 0
$ws2.Cells.Item(4, 58).Value = 0
$ws2.Cells.Item(4, 59).Value = 0
$ws2.Cells.Item(4, 60).Value = 0
$ws2.Cells.Item(4, 61).Value = 0
$ws2.Cells.Item(4, 62).Value = 0
$ws2.Cells.Item(4, 63).Value = 0
$ws2.Cells.Item(4, 64).Value = 0
$ws2.Cells.Item(4, 65).Value = 0
$ws2.Cells.Item(4, 66).Value = 0
$ws2.Cells.Item(4, 67).Value = 0
$ws2.Cells.Item(4, 68).Value = 0
$ws2.Cells.Item(4, 69).Value = 0
$ws2.Cells.Item(4, 70).Value = 0
$ws2.Cells.Item(4, 71).Value = 0
$ws2.Cells.Item(4, 72).Value = 0
$ws2.Cells.Item(4, 73).Value = 0
$ws2.Cells.Item(4, 74).Value = 0
$ws2.Cells.Item(4, 75).Value = 0
$ws2.Cells.Item(4, 76).Value = 0
$ws2.Cells.Item(4, 77).Value = 0
$ws2.Cells.Item(4, 78).Value = 0
$ws2.Cells.Item(4, 79).Value = 0
$ws2.Cells.Item(4, 80).Value = 0
$ws2.Cells.Item(4, 81).Value = 0
$ws2.Cells.Item(4, 82).Value = 0
$ws2.Cells.Item(4, 83).Value = 0
$ws2.Cells.Item(4, 84).Value = 0
$ws2.Cells.Item(4, 85).Value = 0
$ws2.Cells.Item(4, 86).Value = 0
$ws2.Cells.Item(4, 87).Value = 0
$ws2.Cells.Item(4, 88).Value = 0
$ws2.Cells.Item(4, 89).Value = 0
$ws2.Cells.Item(4, 90).Value = 0
$ws2.Cells.Item(4, 91).Value = 0
$ws2.Cells.Item(4, 92).Value = 0
$ws2.Cells.Item(4, 93).Value = 0.1802040037029629
$ws2.Cells.Item(4, 94).Value = 0.1804970386115189
$ws2.Cells.Item(4, 95).Value = 0.3969817138056019
$ws2.Cells.Item(4, 96).Value = 0.4014812945771193
$ws2.Cells.Item(4, 97).Value = 0.4211647509743794
$ws2.Cells.Item(4, 98).Value = 0.4804650757251284
$ws2.Cells.Item(4, 99).Value = 0.4872887329092097
$ws2.Cells.Item(4, 100).Value = 0.4906821964111329
$ws2.Cells.Item(4, 101).Value = 0.4933446761959797
$ws2.Cells.Item(4, 102).Value = 0.5233860273850613
$ws2.Cells.Item(4, 103).Value = 0.5619331155882458
$ws2.Cells.Item(4, 104).Value = 0.585535476282118
$ws2.Cells.Item(4, 105).Value = 0.6581942741104018
$ws2.Cells.Item(4, 106).Value = 0.6598043665834292
$ws2.Cells.Item(4, 107).Value = 0.8166586736750914
$ws2.Cells.Item(4, 108).Value = 0.8453107169342317
$ws2.Cells.Item(4, 109).Value = 0.9199394131446523
$ws2.Cells.Item(4, 110).Value = 0.9432460723993159
$ws2.Cells.Item(4, 111).Value = 0.9534018153403009
$ws2.Cells.Item(4, 112).Value = 0.9534362846942166
$ws2.Cells.Item(4, 113).Value = 0.9536470364541828
$ws2.Cells.Item(4, 114).Value = 0.966504694850306
$ws2.Cells.Item(4, 115).Value = 0.9712171402252134
$ws2.Cells.Item(4, 116).Value = 0.9722122225447842
$ws2.Cells.Item(4, 117).Value = 0.9722374197275118
$ws2.Cells.Item(4, 118).Value = 0.9758590828434176
$ws2.Cells.Item(4, 119).Value = 0.9763792449206188
$ws2.Cells.Item(4, 120).Value = 0.9885250751453114
$ws2.Cells.Item(4, 121).Value = 0.9987667981095857
$ws2.Cells.Item(4, 122).Value = 0.9991294361376469
$ws2.Cells.Item(4, 123).Value = 0.9996599610694884
# row 5
$ws2.Cells.Item(5, 26).Value = 0
$ws2.Cells.Item(5, 27).Value = 0
$ws2.Cells.Item(5, 28).Value = 0
$ws2.Cells.Item(5, 29).Value = 0
$ws2.Cells.Item(5, 30).Value = 0
$ws2.Cells.Item(5, 31).Value = 0
$ws2.Cells.Item(5, 32).Value = 0
$ws2.Cells.Item(5, 33).Value = 0
$ws2.Cells.Item(5, 34).Value = 0
$ws2.Cells.Item(5, 35).Value = 0
$ws2.Cells.Item(5, 36).Value = 0
$ws2.Cells.Item(5, 37).Value = 0
$ws2.Cells.Item(5, 38).Value = 0
$ws2.Cells.Item(5, 39).Value = 0
$ws2.Cells.Item(5, 40).Value = 0
$ws2.Cells.Item(5, 41).Value = 0
$ws2.Cells.Item(5, 42).Value = 0
$ws2.Cells.Item(5, 43).Value = 0
$ws2.Cells.Item(5, 44).Value = 0
$ws2.Cells.Item(5, 45).Value = 0
$ws2.Cells.Item(5, 46).Value = 0
$ws2.Cells.Item(5, 47).Value = 0
$ws2.Cells.Item(5, 48).Value = 0
$ws2.Cells.Item(5, 49).Value = 0
$ws2.Cells.Item(5, 50).Value = 0
$ws2.Cells.Item(5, 51).Value = 0
$ws2.Cells.Item(5, 52).Value = 0
$ws2.Cells.Item(5, 53).Value = 0
$ws2.Cells.Item(5, 54).Value = 0
$ws2.Cells.Item(5, 55).Value = 0
$ws2.Cells.Item(5, 56).Value = 0
$ws2.Cells.Item(5, 57).Value = 0
$ws2.Cells.Item(5, 58).Value = 0
$ws2.Cells.Item(5, 59).Value = 0
$ws2.Cells.Item(5, 60).Value = 0
$ws2.Cells.Item(5, 61).Value = 0
$ws2.Cells.Item(5, 62).Value = 0
$ws2.Cells.Item(5, 63).Value = 0
$ws2.Cells.Item(5, 64).Value = 0
$ws2.Cells.Item(5, 65).Value = 0
$ws2.Cells.Item(5, 66).Value = 0
$ws2.Cells.Item(5, 67).Value = 0
$ws2.Cells.Item(5, 68).Value = 0
$ws2.Cells.Item(5, 69).Value = 0
$ws2.Cells.Item(5, 70).Value = 0
$ws2.Cells.Item(5, 71).Value = 0
$ws2.Cells.Item(5, 72).Value = 0
$ws2.Cells.Item(5, 73).Value = 0
$ws2.Cells.Item(5, 74).Value = 0
$ws2.Cells.Item(5, 75).Value = 0
$ws2.Cells.Item(5, 76).Value = 0
$ws2.Cells.Item(5, 77).Value = 0
$ws2.Cells.Item(5, 78).Value = 0
$ws2.Cells.Item(5, 79).Value = 0
$ws2.Cells.Item(5, 80).Value = 0
$ws2.Cells.Item(5, 81).Value = 0
$ws2.Cells.Item(5, 82).Value = 0
$ws2.Cells.Item(5, 83).Value = 0
$ws2.Cells.Item(5, 84).Value = 0
$ws2.Cells.Item(5, 85).Value = 0
$ws2.Cells.Item(5, 86).Value = 0
$ws2.Cells.Item(5, 87).Value = 0
$ws2.Cells.Item(5, 88).Value = 0
$ws2.Cells.Item(5, 89).Value = 0
$ws2.Cells.Item(5, 90).Value = 0
$ws2.Cells.Item(5, 91).Value = 0
$ws2.Cells.Item(5, 92).Value = 0
$ws2.Cells.Item(5, 93).Value = 0.0414397235972116
$ws2.Cells.Item(5, 94).Value = 0.3045203506098182
$ws2.Cells.Item(5, 95).Value = 0.5058136796787165
$ws2.Cells.Item(5, 96).Value = 0.5594627751850137
$ws2.Cells.Item(5, 97).Value = 0.5606530135377783
$ws2.Cells.Item(5, 98).Value = 0.6043563704923376
$ws2.Cells.Item(5, 99).Value = 0.6246117891629708
$ws2.Cells.Item(5, 100).Value = 0.66973085520673
$ws2.Cells.Item(5, 101).Value = 0.6827649700748417
$ws2.Cells.Item(5, 102).Value = 0.705196277165789
$ws2.Cells.Item(5, 103).Value = 0.7476651622896823
$ws2.Cells.Item(5, 104).Value = 0.7482238356174075
$ws2.Cells.Item(5, 105).Value = 0.818077776555288
$ws2.Cells.Item(5, 106).Value = 0.8190547298348462
$ws2.Cells.Item(5, 107).Value = 0.9306554009696474
$ws2.Cells.Item(5, 108).Value = 0.9395786715785355
$ws2.Cells.Item(5, 109).Value = 0.9520444328320571
$ws2.Cells.Item(5, 110).Value = 0.9591281254185006
$ws2.Cells.Item(5, 111).Value = 0.9634073498243393
$ws2.Cells.Item(5, 112).Value = 0.9639880451272503
$ws2.Cells.Item(5, 113).Value = 0.9767439497670768
$ws2.Cells.Item(5, 114).Value = 0.9789881342592095
$ws2.Cells.Item(5, 115).Value = 0.9793661320647564
$ws2.Cells.Item(5, 116).Value = 0.9795267487165394
$ws2.Cells.Item(5, 117).Value = 0.982852191122305
$ws2.Cells.Item(5, 118).Value = 0.9828545130543153
$ws2.Cells.Item(5, 119).Value = 0.9840562400201567
$ws2.Cells.Item(5, 120).Value = 0.9867973824562909
$ws2.Cells.Item(5, 121).Value = 0.9867986577606321
$ws2.Cells.Item(5, 122).Value = 0.987480236076559
$ws2.Cells.Item(5, 123).Value = 0.9993505571765627
$ws2.Cells.Item(5, 124).Value = 1.0
$ws2.Cells.Item(5, 125).Value = 1.0
$ws2.Cells.Item(5, 126).Value = 1.0
# row 8
$ws2.Cells.Item(8, 17).Value = 0
$ws2.Cells.Item(8, 18).Value = 0
$ws2.Cells.Item(8, 19).Value = 0
$ws2.Cells.Item(8, 20).Value = 0
$ws2.Cells.Item(8, 21).Value = 0
$ws2.Cells.Item(8, 22).Value = 0
$ws2.Cells.Item(8, 23).Value = 0
$ws2.Cells.Item(8, 24).Value = 0
$ws2.Cells.Item(8, 25).Value = 0
$ws2.Cells.Item(8, 26).Value = 0
$ws2.Cells.Item(8, 27).Value = 0
$ws2.Cells.Item(8, 28).Value = 0
$ws2.Cells.Item(8, 29).Value = 0
$ws2.Cells.Item(8, 30).Value = 0
$ws2.Cells.Item(8, 31).Value = 0
$ws2.Cells.Item(8, 32).Value = 0
$ws2.Cells.Item(8, 33).Value = 0
$ws2.Cells.Item(8, 34).Value = 0
$ws2.Cells.Item(8, 35).Value = 0
$ws2.Cells.Item(8, 36).Value = 0
$ws2.Cells.Item(8, 37).Value = 0
$ws2.Cells.Item(8, 38).Value = 0
$ws2.Cells.Item(8, 39).Value = 0
$ws2.Cells.Item(8, 40).Value = 0
$ws2.Cells.Item(8, 41).Value = 0
$ws2.Cells.Item(8, 42).Value = 0
$ws2.Cells.Item(8, 43).Value = 0
$ws2.Cells.Item(8, 44).Value = 0
$ws2.Cells.Item(8, 45).Value = 0
$ws2.Cells.Item(8, 46).Value = 0
$ws2.Cells.Item(8, 47).Value = 0
$ws2.Cells.Item(8, 48).Value = 0
$ws2.Cells.Item(8, 49).Value = 0
$ws2.Cells.Item(8, 50).Value = 0
$ws2.Cells.Item(8, 51).Value = 0
$ws2.Cells.Item(8, 52).Value = 0
$ws2.Cells.Item(8, 53).Value = 0
$ws2.Cells.Item(8, 54).Value = 0
$ws2.Cells.Item(8, 55).Value = 0
$ws2.Cells.Item(8, 56).Value = 0
$ws2.Cells.Item(8, 57).Value = 0
$ws2.Cells.Item(8, 58).Value = 0
$ws2.Cells.Item(8, 59).Value = 0
$ws2.Cells.Item(8, 60).Value = 0
$ws2.Cells.Item(8, 61).Value = 0
$ws2.Cells.Item(8, 62).Value = 0
$ws2.Cells.Item(8, 63).Value = 0
$ws2.Cells.Item(8, 64).Value = 0
$ws2.Cells.Item(8, 65).Value = 0
$ws2.Cells.Item(8, 66).Value = 0
$ws2.Cells.Item(8, 67).Value = 0
$ws2.Cells.Item(8, 68).Value = 0
$ws2.Cells.Item(8, 69).Value = 0
$ws2.Cells.Item(8, 70).Value = 0
$ws2.Cells.Item(8, 71).Value = 0
$ws2.Cells.Item(8, 72).Value = 0
$ws2.Cells.Item(8, 73).Value = 0
$ws2.Cells.Item(8, 74).Value = 0
$ws2.Cells.Item(8, 75).Value = 0
$ws2.Cells.Item(8, 76).Value = 0
$ws2.Cells.Item(8, 77).Value = 0
$ws2.Cells.Item(8, 78).Value = 0
$ws2.Cells.Item(8, 79).Value = 0
$ws2.Cells.Item(8, 80).Value = 0
$ws2.Cells.Item(8, 81).Value = 0
$ws2.Cells.Item(8, 82).Value = 0
$ws2.Cells.Item(8, 83).Value = 0
$ws2.Cells.Item(8, 84).Value = 0
$ws2.Cells.Item(8, 85).Value = 0
$ws2.Cells.Item(8, 86).Value = 0
$ws2.Cells.Item(8, 87).Value = 0
$ws2.Cells.Item(8, 88).Value = 0
$ws2.Cells.Item(8, 89).Value = 0
$ws2.Cells.Item(8, 90).Value = 0
$ws2.Cells.Item(8, 91).Value = 0
$ws2.Cells.Item(8, 92).Value = 0.14584622503048
$ws2.Cells.Item(8, 93).Value = 0.2499583598719811
$ws2.Cells.Item(8, 94).Value = 0.5100051977118101
$ws2.Cells.Item(8, 95).Value = 0.5519824553787543
$ws2.Cells.Item(8, 96).Value = 0.5525707117810629
$ws2.Cells.Item(8, 97).Value = 0.581803084852136
$ws2.Cells.Item(8, 98).Value = 0.604459840125297
$ws2.Cells.Item(8, 99).Value = 0.6188393558912532
$ws2.Cells.Item(8, 100).Value = 0.6194324202052556
$ws2.Cells.Item(8, 101).Value = 0.7418526027076126
$ws2.Cells.Item(8, 102).Value = 0.7803667392499316
$ws2.Cells.Item(8, 103).Value = 0.7999827393659655
$ws2.Cells.Item(8, 104).Value = 0.8483565804630315
$ws2.Cells.Item(8, 105).Value = 0.8488448079511451
$ws2.Cells.Item(8, 106).Value = 0.9104203347053635
$ws2.Cells.Item(8, 107).Value = 0.9104350369777624
$ws2.Cells.Item(8, 108).Value = 0.920681983308746
$ws2.Cells.Item(8, 109).Value = 0.9270307966198672
$ws2.Cells.Item(8, 110).Value = 0.9316444924174939
$ws2.Cells.Item(8, 111).Value = 0.9367878577262195
$ws2.Cells.Item(8, 112).Value = 0.9491935903547631
$ws2.Cells.Item(8, 113).Value = 0.9605656126136474
$ws2.Cells.Item(8, 114).Value = 0.9636167408529507
$ws2.Cells.Item(8, 115).Value = 0.965548790370621
$ws2.Cells.Item(8, 116).Value = 0.9660341320692512
$ws2.Cells.Item(8, 117).Value = 0.9747417992894751
$ws2.Cells.Item(8, 118).Value = 0.9891175899730182
$ws2.Cells.Item(8, 119).Value = 0.9976371041786777
$ws2.Cells.Item(8, 120).Value = 0.9983399393955109
$ws2.Cells.Item(8, 121).Value = 0.9985383164996472
$ws2.Cells.Item(8, 122).Value = 0.9993203757406766
# row 11
$ws2.Cells.Item(11, 20).Value = 0
$ws2.Cells.Item(11, 21).Value = 0
$ws2.Cells.Item(11, 22).Value = 0
$ws2.Cells.Item(11, 23).Value = 0
$ws2.Cells.Item(11, 24).Value = 0
$ws2.Cells.Item(11, 25).Value = 0
$ws2.Cells.Item(11, 26).Value = 0
$ws2.Cells.Item(11, 27).Value = 0
$ws2.Cells.Item(11, 28).Value = 0
$ws2.Cells.Item(11, 29).Value = 0
$ws2.Cells.Item(11, 30).Value = 0
$ws2.Cells.Item(11, 31).Value = 0
$ws2.Cells.Item(11, 32).Value = 0
$ws2.Cells.Item(11, 33).Value = 0
$ws2.Cells.Item(11, 34).Value = 0
$ws2.Cells.Item(11, 35).Value = 0
$ws2.Cells.Item(11, 36).Value = 0
$ws2.Cells.Item(11, 37).Value = 0
$ws2.Cells.Item(11, 38).Value = 0
$ws2.Cells.Item(11, 39).Value = 0
$ws2.Cells.Item(11, 40).Value = 0
$ws2.Cells.Item(11, 41).Value = 0
$ws2.Cells.Item(11, 42).Value = 0
$ws2.Cells.Item(11, 43).Value = 0
$ws2.Cells.Item(11, 44).Value = 0
$ws2.Cells.Item(11, 45).Value = 0
$ws2.Cells.Item(11, 46).Value = 0
$ws2.Cells.Item(11, 47).Value = 0
$ws2.Cells.Item(11, 48).Value = 0
$ws2.Cells.Item(11, 49).Value = 0
$ws2.Cells.Item(11, 50).Value = 0
$ws2.Cells.Item(11, 51).Value = 0
$ws2.Cells.Item(11, 52).Value = 0.2230194153683032
$ws2.Cells.Item(11, 53).Value = 0.2260492096176228
$ws2.Cells.Item(11, 54).Value = 0.5238944720853328
$ws2.Cells.Item(11, 55).Value = 0.5671828455662149
$ws2.Cells.Item(11, 56).Value = 0.5744841623486575
$ws2.Cells.Item(11, 57).Value = 0.651159379344117
$ws2.Cells.Item(11, 58).Value = 0.6772736197498181
$ws2.Cells.Item(11, 59).Value = 0.6776749489258015
$ws2.Cells.Item(11, 60).Value = 0.6813062560630447
$ws2.Cells.Item(11, 61).Value = 0.794326226728748
$ws2.Cells.Item(11, 62).Value = 0.8079013175217049
$ws2.Cells.Item(11, 63).Value = 0.8361115555461832
$ws2.Cells.Item(11, 64).Value = 0.8802093134309475
$ws2.Cells.Item(11, 65).Value = 0.880252348264749
$ws2.Cells.Item(11, 66).Value = 0.9087535940458672
$ws2.Cells.Item(11, 67).Value = 0.9150865407218038
$ws2.Cells.Item(11, 68).Value = 0.9176856713718896
$ws2.Cells.Item(11, 69).Value = 0.9184754118097952
$ws2.Cells.Item(11, 70).Value = 0.9193325099381688
$ws2.Cells.Item(11, 71).Value = 0.9193646063115676
$ws2.Cells.Item(11, 72).Value = 0.9227566657824392
$ws2.Cells.Item(11, 73).Value = 0.9259866581803353
$ws2.Cells.Item(11, 74).Value = 0.9263662821091108
$ws2.Cells.Item(11, 75).Value = 0.9555306448657819
$ws2.Cells.Item(11, 76).Value = 0.9579684834172699
$ws2.Cells.Item(11, 77).Value = 0.9660781949195535
$ws2.Cells.Item(11, 78).Value = 0.9841574011104945
$ws2.Cells.Item(11, 79).Value = 0.9891420284347742
$ws2.Cells.Item(11, 80).Value = 0.9891422605126011
$ws2.Cells.Item(11, 81).Value = 0.9891978348785391
$ws2.Cells.Item(11, 82).Value = 0.9913873437331517

# --- Step3_DataPts_* sheets: update summary columns C, D, F, G (rows 4,5,8,11) ---
$ws3 = $wb.Worksheets.Item("Step3_DataPts_0.5")
# row 4
$ws3.Cells.Item(4, 3).Value = 90.0
$ws3.Cells.Item(4, 4).Value = 101.0
$ws3.Cells.Item(4, 6).Value = 0.5233860273850613
$ws3.Cells.Item(4, 7).Value = 11.0
# row 5
$ws3.Cells.Item(5, 3).Value = 91.0
$ws3.Cells.Item(5, 4).Value = 94.0
$ws3.Cells.Item(5, 6).Value = 0.5058136796787165
$ws3.Cells.Item(5, 7).Value = 3.0
# row 8
$ws3.Cells.Item(8, 3).Value = 89.0
$ws3.Cells.Item(8, 4).Value = 93.0
$ws3.Cells.Item(8, 6).Value = 0.5100051977118101
$ws3.Cells.Item(8, 7).Value = 4.0
# row 11
$ws3.Cells.Item(11, 3).Value = 49.0
$ws3.Cells.Item(11, 4).Value = 53.0
$ws3.Cells.Item(11, 6).Value = 0.5238944720853328
$ws3.Cells.Item(11, 7).Value = 4.0

$ws4 = $wb.Worksheets.Item("Step3_DataPts_0.7")
# row 4
$ws4.Cells.Item(4, 3).Value = 90.0
$ws4.Cells.Item(4, 4).Value = 106.0
$ws4.Cells.Item(4, 6).Value = 0.8166586736750914
$ws4.Cells.Item(4, 7).Value = 16.0
# row 5
$ws4.Cells.Item(5, 3).Value = 91.0
$ws4.Cells.Item(5, 4).Value = 101.0
$ws4.Cells.Item(5, 6).Value = 0.705196277165789
$ws4.Cells.Item(5, 7).Value = 10.0
# row 8
$ws4.Cells.Item(8, 3).Value = 89.0
$ws4.Cells.Item(8, 4).Value = 100.0
$ws4.Cells.Item(8, 6).Value = 0.7418526027076126
$ws4.Cells.Item(8, 7).Value = 11.0
# row 11
$ws4.Cells.Item(11, 3).Value = 49.0
$ws4.Cells.Item(11, 4).Value = 60.0
$ws4.Cells.Item(11, 6).Value = 0.794326226728748
$ws4.Cells.Item(11, 7).Value = 11.0

$ws5 = $wb.Worksheets.Item("Step3_DataPts_0.8")
# row 4
$ws5.Cells.Item(4, 3).Value = 90.0
$ws5.Cells.Item(4, 4).Value = 106.0
$ws5.Cells.Item(4, 6).Value = 0.8166586736750914
$ws5.Cells.Item(4, 7).Value = 16.0
# row 5
$ws5.Cells.Item(5, 3).Value = 91.0
$ws5.Cells.Item(5, 4).Value = 104.0
$ws5.Cells.Item(5, 6).Value = 0.818077776555288
$ws5.Cells.Item(5, 7).Value = 13.0
# row 8
$ws5.Cells.Item(8, 3).Value = 89.0
$ws5.Cells.Item(8, 4).Value = 103.0
$ws5.Cells.Item(8, 6).Value = 0.8483565804630315
$ws5.Cells.Item(8, 7).Value = 14.0
# row 11
$ws5.Cells.Item(11, 3).Value = 49.0
$ws5.Cells.Item(11, 4).Value = 61.0
$ws5.Cells.Item(11, 6).Value = 0.8079013175217049
$ws5.Cells.Item(11, 7).Value = 12.0

$ws6 = $wb.Worksheets.Item("Step3_DataPts_0.9")
# row 4
$ws6.Cells.Item(4, 3).Value = 90.0
$ws6.Cells.Item(4, 4).Value = 108.0
$ws6.Cells.Item(4, 6).Value = 0.9199394131446523
$ws6.Cells.Item(4, 7).Value = 18.0
# row 5
$ws6.Cells.Item(5, 3).Value = 91.0
$ws6.Cells.Item(5, 4).Value = 106.0
$ws6.Cells.Item(5, 6).Value = 0.9306554009696474
$ws6.Cells.Item(5, 7).Value = 15.0
# row 8
$ws6.Cells.Item(8, 3).Value = 89.0
$ws6.Cells.Item(8, 4).Value = 105.0
$ws6.Cells.Item(8, 6).Value = 0.9104203347053635
$ws6.Cells.Item(8, 7).Value = 16.0
# row 11
$ws6.Cells.Item(11, 3).Value = 49.0
$ws6.Cells.Item(11, 4).Value = 65.0
$ws6.Cells.Item(11, 6).Value = 0.9087535940458672
$ws6.Cells.Item(11, 7).Value = 16.0

